$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns O:P, matching the style used by the existing header row ---
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("O1").Value = "MensajeUsuario"
$ws.Range("P1").Value = "RespuestaIA"

# --- Existing rows 2-4 get blank cells in the two new columns ---
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = ""
$ws.Range("O3").Value = ""
$ws.Range("P3").Value = ""
$ws.Range("O4").Value = ""
$ws.Range("P4").Value = ""

# --- B4 becomes a true numeric value instead of text ---
$ws.Range("B4").Value = 1000271912

# --- New row 5: latest negotiation log entry ---
$ws.Range("A5").Value = "2025-10-15 22:19:16"
# Cedula is kept as text (quote-prefixed) like the original rows, not auto-converted to a number
$ws.Range("B5").Value = "'1000274330"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "Elian"
$ws.Range("D5").Value = "TARJETA DE CRÉDITO"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "PRORROGA SIN PAGO"
$ws.Range("G5").Value = "48 cuotas"
$ws.Range("H5").Value = "34.19.100.134"
$ws.Range("I5").Value = "The Dalles"
$ws.Range("J5").Value = "Oregon"
$ws.Range("K5").Value = "United States"
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = ""
$ws.Range("P5").Value = ""
